# AlojamientoSeleccionado.xlsx - swap the sample listing shown in the
# "Alojamiento seleccionado" confirmation sheet so the automated test no
# longer depends on a specific (changeable) Airbnb-style listing: update
# the apartment name, host name and price, and shrink column A now that
# the new name is a bit shorter than the old "bestFit" width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Apartamento tipo Ático, piso 4 SIN Ascensor"
$ws.Range("B2").Value = "Hey You"
$ws.Range("C2").Value = "$1,022,861.32 COP"

# Column A's bestFit width shrinks to ~40.56 chars for the new text.
# (Re-assert B/C's existing widths too so they keep their explicit
# custom-width flag instead of reverting to "not custom" after the edit.)
$ws.Columns.Item(1).ColumnWidth = 39.666666666666664
$ws.Columns.Item(2).ColumnWidth = 8.166666666666668
$ws.Columns.Item(3).ColumnWidth = 17.333333333333336
